$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Sheet "CypherOutput": replace the existing case rows (2-6) with the new
# canine-filter data set, and append a new row (7) for the Glioma case.
# All cells in this table are text (e.g. ages like "11.5" are strings, not
# numbers), so force text formatting before writing the values.
# ---------------------------------------------------------------------------
$wsCypher = $wb.Worksheets.Item("CypherOutput")

# Only the Age column (G) holds numeric-looking text (e.g. "11.5"); force it
# to text so it is not auto-converted to a number. The other columns contain
# non-numeric text and do not need special formatting.
$wsCypher.Range("G2:G7").NumberFormat = "@"

# Row 2: COTC007B-0503
$wsCypher.Range("A2").Value = "COTC007B-0503"
$wsCypher.Range("B2").Value = "COTC007B"
$wsCypher.Range("C2").Value = "Clinical Trial"
$wsCypher.Range("D2").Value = "Beagle"
$wsCypher.Range("E2").Value = "Lymphoma"
$wsCypher.Range("F2").Value = "IIIa"
$wsCypher.Range("G2").Value = "11.5"
$wsCypher.Range("H2").Value = "Female"
$wsCypher.Range("I2").Value = "Yes"

# Row 3: COTC007B-0211
$wsCypher.Range("A3").Value = "COTC007B-0211"
$wsCypher.Range("B3").Value = "COTC007B"
$wsCypher.Range("C3").Value = "Clinical Trial"
$wsCypher.Range("D3").Value = "Beagle"
$wsCypher.Range("E3").Value = "Lymphoma"
$wsCypher.Range("F3").Value = "III"
$wsCypher.Range("G3").Value = "9.9"
$wsCypher.Range("H3").Value = "Male"
$wsCypher.Range("I3").Value = "Yes"

# Row 4: COTC007B-0510
$wsCypher.Range("A4").Value = "COTC007B-0510"
$wsCypher.Range("B4").Value = "COTC007B"
$wsCypher.Range("C4").Value = "Clinical Trial"
$wsCypher.Range("D4").Value = "Beagle"
$wsCypher.Range("E4").Value = "Lymphoma"
$wsCypher.Range("F4").Value = "IIIa"
$wsCypher.Range("G4").Value = "9.5"
$wsCypher.Range("H4").Value = "Male"
$wsCypher.Range("I4").Value = "Yes"

# Row 5: COTC007B-0608
$wsCypher.Range("A5").Value = "COTC007B-0608"
$wsCypher.Range("B5").Value = "COTC007B"
$wsCypher.Range("C5").Value = "Clinical Trial"
$wsCypher.Range("D5").Value = "Beagle"
$wsCypher.Range("E5").Value = "Lymphoma"
$wsCypher.Range("F5").Value = "IVa"
$wsCypher.Range("G5").Value = "4.2"
$wsCypher.Range("H5").Value = "Male"
$wsCypher.Range("I5").Value = "Yes"

# Row 6: NCATS-COP01-CCB010045
$wsCypher.Range("A6").Value = "NCATS-COP01-CCB010045"
$wsCypher.Range("B6").Value = "NCATS-COP01"
$wsCypher.Range("C6").Value = "Transcriptomics"
$wsCypher.Range("D6").Value = "Beagle"
$wsCypher.Range("E6").Value = "T Cell Lymphoma"
$wsCypher.Range("F6").Value = "III"
$wsCypher.Range("G6").Value = "5.0"
$wsCypher.Range("H6").Value = "Female"
$wsCypher.Range("I6").Value = "Yes"

# Row 7 (new): GLIOMA01-i_6561
$wsCypher.Range("A7").Value = "GLIOMA01-i_6561"
$wsCypher.Range("B7").Value = "GLIOMA01"
$wsCypher.Range("C7").Value = "Genomics"
$wsCypher.Range("D7").Value = "Beagle"
$wsCypher.Range("E7").Value = "Glioma"
$wsCypher.Range("F7").Value = "Unknown"
$wsCypher.Range("G7").Value = "8.0"
$wsCypher.Range("H7").Value = "Male"
$wsCypher.Range("I7").Value = "Yes"

# ---------------------------------------------------------------------------
# Sheet "StatOutput": update the summary counts (also stored as text).
# ---------------------------------------------------------------------------
$wsStat = $wb.Worksheets.Item("StatOutput")

$wsStat.Range("A2:D2").NumberFormat = "@"
$wsStat.Range("A2").Value = "12"
$wsStat.Range("B2").Value = "8"
$wsStat.Range("C2").Value = "6"
$wsStat.Range("D2").Value = "3"

# ---------------------------------------------------------------------------
# Sheet "CaseDetailStat": the header row (file metadata column titles) is
# cleared out, leaving the sheet with only the data row.
# ---------------------------------------------------------------------------
$wsCaseDetail = $wb.Worksheets.Item("CaseDetailStat")
$wsCaseDetail.Range("A1:F1").ClearContents()

# ---------------------------------------------------------------------------
# Sheet "CaseDetailStat_Message": a new error line is inserted ahead of the
# third Neo4j connection/query message block, and the Cypher query value in
# that block is now blank.
# ---------------------------------------------------------------------------
$wsCaseDetailMsg = $wb.Worksheets.Item("CaseDetailStat_Message")
$wsCaseDetailMsg.Rows(21).Insert()
$wsCaseDetailMsg.Range("A21").Value = "Cypher query should not be an empty string"
$wsCaseDetailMsg.Range("A29").Value = ""
